$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 7694241.5
$ws.Range("I19").Value = 1653
$ws.Range("K19").Value = 1653
$ws.Range("M19").Value = -1478
$ws.Range("H33").Value = 284.7
$ws.Range("I33").Value = 258.35715
$ws.Range("K33").Value = 258.35715
$ws.Range("M33").Value = -29.35714999999999
$ws.Range("H64").Value = 4041.4443
$ws.Range("I64").Value = 3934.5
$ws.Range("J64").Value = 4072
$ws.Range("K64").Value = 3934.5
$ws.Range("L64").Value = 4072
$ws.Range("M64").Value = -3686.5
$ws.Range("N64").Value = -4568
$ws.Range("H67").Value = 4041.4443
$ws.Range("I67").Value = 3934.5
$ws.Range("J67").Value = 4072
$ws.Range("K67").Value = 3934.5
$ws.Range("L67").Value = 4072
$ws.Range("M67").Value = -3076.5
$ws.Range("N67").Value = -5788
$ws.Range("H69").Value = 16221.889
$ws.Range("I69").Value = 0
$ws.Range("J69").Value = 16221.889
$ws.Range("K69").Value = 0
$ws.Range("L69").Value = 48665.667
$ws.Range("N69").Value = -50413.667
$ws.Range("H72").Value = 16221.889
$ws.Range("I72").Value = 0
$ws.Range("J72").Value = 16221.889
$ws.Range("K72").Value = 0
$ws.Range("L72").Value = 145997.001
$ws.Range("N72").Value = -154733.001
$ws.Range("H88").Value = 2752.7334
$ws.Range("J88").Value = 2620.7778
$ws.Range("L88").Value = 2620.7778
$ws.Range("N88").Value = -3432.7778
$ws.Range("H91").Value = 2752.7334
$ws.Range("J91").Value = 2620.7778
$ws.Range("L91").Value = 2620.7778
$ws.Range("N91").Value = -5428.7778
$ws.Range("H94").Value = 2906.625
$ws.Range("I94").Value = 2906.625
$ws.Range("K94").Value = 2906.625
$ws.Range("M94").Value = -2455.625
$ws.Range("H100").Value = 7399.6
$ws.Range("J100").Value = 8999.333000000001
$ws.Range("L100").Value = 8999.333000000001
$ws.Range("N100").Value = -10081.333
$ws.Range("H129").Value = 2998.8147
$ws.Range("I129").Value = 775.75
$ws.Range("J129").Value = 6232.364
$ws.Range("K129").Value = 2327.25
$ws.Range("L129").Value = 18697.092
$ws.Range("M129").Value = 2672.75
$ws.Range("N129").Value = -28697.092
$ws.Range("H137").Value = 1944.1052
$ws.Range("I137").Value = 1906.9412
$ws.Range("J137").Value = 2260
$ws.Range("K137").Value = 5720.8236
$ws.Range("L137").Value = 6780
$ws.Range("M137").Value = -3170.8236
$ws.Range("N137").Value = -11880
$ws.Range("H141").Value = 17881030
$ws.Range("I141").Value = 25006894
$ws.Range("K141").Value = 75020682
$ws.Range("M141").Value = -75015502
$ws.Range("M69").ClearContents()
$ws.Range("M72").ClearContents()

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3125.348
$ws.Range("I2").Value = 3090.2144
$ws.Range("J2").Value = 3180
$ws.Range("K2").Value = 3090.2144
$ws.Range("L2").Value = 3180
$ws.Range("M2").Value = -2977.2144
$ws.Range("N2").Value = -3406
$ws.Range("H30").Value = 40250
$ws.Range("J30").Value = 80000
$ws.Range("L30").Value = 80000
$ws.Range("N30").Value = -80300
$ws.Range("H32").Value = 3168.5076
$ws.Range("I32").Value = 2323.61
$ws.Range("K32").Value = 2323.61
$ws.Range("M32").Value = -2036.61
$ws.Range("H45").Value = 8000
$ws.Range("I45").Value = 0
$ws.Range("K45").Value = 0
$ws.Range("H97").Value = 1761.8125
$ws.Range("I97").Value = 1862.2307
$ws.Range("J97").Value = 1326.6666
$ws.Range("K97").Value = 1862.2307
$ws.Range("L97").Value = 1326.6666
$ws.Range("M97").Value = -1366.2307
$ws.Range("N97").Value = -2318.6666
$ws.Range("H116").Value = 3125.348
$ws.Range("I116").Value = 3090.2144
$ws.Range("J116").Value = 3180
$ws.Range("K116").Value = 3090.2144
$ws.Range("L116").Value = 3180
$ws.Range("M116").Value = -796.2143999999998
$ws.Range("N116").Value = -7768
$ws.Range("H132").Value = 2385962.5
$ws.Range("I132").Value = 4548.4116
$ws.Range("K132").Value = 13645.2348
$ws.Range("M132").Value = -11115.2348
$ws.Range("M45").ClearContents()

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3125.348
$ws.Range("I3").Value = 3090.2144
$ws.Range("J3").Value = 3180
$ws.Range("K3").Value = 3090.2144
$ws.Range("L3").Value = 3180
$ws.Range("M3").Value = -2976.2144
$ws.Range("N3").Value = -3408
$ws.Range("H105").Value = 717748.9
$ws.Range("I105").Value = 1429464.9
$ws.Range("K105").Value = 1429464.9
$ws.Range("M105").Value = -1427717.9

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 706.8333
$ws.Range("I7").Value = 663.3333
$ws.Range("K7").Value = 663.3333
$ws.Range("M7").Value = -550.3333
$ws.Range("H31").Value = 38465136
$ws.Range("I31").Value = 66670212
$ws.Range("K31").Value = 66670212
$ws.Range("M31").Value = -66669917
$ws.Range("H34").Value = 38465136
$ws.Range("I34").Value = 66670212
$ws.Range("K34").Value = 66670212
$ws.Range("M34").Value = -66670010
$ws.Range("H39").Value = 0
$ws.Range("I39").Value = 0
$ws.Range("J39").Value = 0
$ws.Range("K39").Value = 0
$ws.Range("L39").Value = 0
$ws.Range("H42").Value = 92000
$ws.Range("J42").Value = 92000
$ws.Range("L42").Value = 92000
$ws.Range("N42").Value = -93186
$ws.Range("H49").Value = 0
$ws.Range("I49").Value = 0
$ws.Range("J49").Value = 0
$ws.Range("K49").Value = 0
$ws.Range("L49").Value = 0
$ws.Range("H99").Value = 40156.223
$ws.Range("I99").Value = 10501.2
$ws.Range("J99").Value = 77225
$ws.Range("K99").Value = 10501.2
$ws.Range("L99").Value = 77225
$ws.Range("M99").Value = -9003.200000000001
$ws.Range("N99").Value = -80221
$ws.Range("H126").Value = 40156.223
$ws.Range("I126").Value = 10501.2
$ws.Range("J126").Value = 77225
$ws.Range("K126").Value = 31503.6
$ws.Range("L126").Value = 231675
$ws.Range("M126").Value = -29033.6
$ws.Range("N126").Value = -236615
$ws.Range("M39").ClearContents()
$ws.Range("N39").ClearContents()
$ws.Range("M49").ClearContents()
$ws.Range("N49").ClearContents()

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 3055.625
$ws.Range("I5").Value = 2800.2
$ws.Range("J5").Value = 3481.3333
$ws.Range("K5").Value = 8400.599999999999
$ws.Range("L5").Value = 10443.9999
$ws.Range("M5").Value = -8288.599999999999
$ws.Range("N5").Value = -10667.9999
$ws.Range("H131").Value = 4109.6
$ws.Range("I131").Value = 2238.4614
$ws.Range("J131").Value = 7584.5713
$ws.Range("K131").Value = 6715.3842
$ws.Range("L131").Value = 22753.7139
$ws.Range("M131").Value = -1675.3842
$ws.Range("N131").Value = -32833.7139
$ws.Range("H135").Value = 3055.625
$ws.Range("I135").Value = 2800.2
$ws.Range("J135").Value = 3481.3333
$ws.Range("K135").Value = 25201.8
$ws.Range("L135").Value = 31331.9997
$ws.Range("M135").Value = -22666.8
$ws.Range("N135").Value = -36401.9997
$ws.Range("H138").Value = 22906.166
$ws.Range("J138").Value = 17176.5
$ws.Range("L138").Value = 51529.5
$ws.Range("N138").Value = -61809.5
$ws.Range("H140").Value = 4323.909
$ws.Range("I140").Value = 1423
$ws.Range("K140").Value = 4269
$ws.Range("M140").Value = 911

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2638.625
$ws.Range("I102").Value = 2590.9644
$ws.Range("K102").Value = 2590.9644
$ws.Range("M102").Value = -968.9643999999998

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2875.4138
$ws.Range("I40").Value = 2905.36
$ws.Range("K40").Value = 2905.36
$ws.Range("M40").Value = -2769.36
$ws.Range("H68").Value = 10418917
$ws.Range("J68").Value = 4003
$ws.Range("L68").Value = 4003
$ws.Range("N68").Value = -5501
$ws.Range("H71").Value = 10418917
$ws.Range("J71").Value = 4003
$ws.Range("L71").Value = 20015
$ws.Range("N71").Value = -27503
$ws.Range("H82").Value = 4353.75
$ws.Range("I82").Value = 3352.4285
$ws.Range("J82").Value = 5132.5557
$ws.Range("K82").Value = 3352.4285
$ws.Range("L82").Value = 5132.5557
$ws.Range("M82").Value = -2991.4285
$ws.Range("N82").Value = -5854.5557
$ws.Range("H85").Value = 4353.75
$ws.Range("I85").Value = 3352.4285
$ws.Range("J85").Value = 5132.5557
$ws.Range("K85").Value = 3352.4285
$ws.Range("L85").Value = 5132.5557
$ws.Range("M85").Value = -2104.4285
$ws.Range("N85").Value = -7628.5557
$ws.Range("H93").Value = 6951807
$ws.Range("I93").Value = 4600.5
$ws.Range("K93").Value = 4600.5
$ws.Range("M93").Value = -3352.5
$ws.Range("H100").Value = 35753720
$ws.Range("I100").Value = 4024.5
$ws.Range("J100").Value = 83419980
$ws.Range("K100").Value = 4024.5
$ws.Range("L100").Value = 83419980
$ws.Range("M100").Value = -3483.5
$ws.Range("N100").Value = -83421062
$ws.Range("H122").Value = 3620.5
$ws.Range("I122").Value = 3315.5405
$ws.Range("K122").Value = 9946.621500000001
$ws.Range("M122").Value = -7496.621500000001

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 21493.666
$ws.Range("J41").Value = 20393.4
$ws.Range("L41").Value = 20393.4
$ws.Range("N41").Value = -21173.4
$ws.Range("H122").Value = 2276.9
$ws.Range("I122").Value = 2252.1667
$ws.Range("J122").Value = 2499.5
$ws.Range("K122").Value = 6756.500100000001
$ws.Range("L122").Value = 7498.5
$ws.Range("M122").Value = -4306.500100000001
$ws.Range("N122").Value = -12398.5
$ws.Range("H132").Value = 284911.5
$ws.Range("I132").Value = 6897.387
$ws.Range("K132").Value = 20692.161
$ws.Range("M132").Value = -18162.161
